$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value
$updates = [ordered]@{
    "D2" = "68.623.52"
    "E2" = "  +1.24%  "
    "D3" = "3.813.64"
    "E3" = "  -0.10%  "
    "E4" = "  +0.06%  "
    "D5" = "609.02"
    "E5" = "  +0.83%  "
    "D6" = "164.54"
    "E6" = "  -1.12%  "
    "D7" = "3.811.98"
    "E7" = "  -0.02%  "
    "E8" = "  -0.05%  "
    "E9" = "  -0.20%  "
    "D10" = "0.160"
    "E10" = "  -0.09%  "
    "D11" = "7.00"
    "E11" = "  +11.26%  "
    "D12" = "0.452"
    "E12" = "  -0.07%  "
    "D13" = "0.0000249"
    "E13" = "  -1.35%  "
    "D14" = "35.25"
    "E14" = "  -2.05%  "
    "D15" = "4.454.82"
    "E15" = "  -0.14%  "
    "D16" = "3.781.31"
    "E16" = "  -1.24%  "
    "D17" = "68.600.33"
    "E17" = "  +1.16%  "
    "D18" = "18.16"
    "E18" = "  -1.75%  "
    "E19" = "  +1.98%  "
    "D20" = "7.08"
    "E20" = "  +0.10%  "
    "D21" = "464.33"
    "E21" = "  +0.35%  "
    "D22" = "9.63"
    "E22" = "  -2.55%  "
    "D23" = "0.699"
    "E23" = "  -0.40%  "
    "D24" = "0.0000149"
    "E24" = "  +0.72%  "
    "D25" = "83.78"
    "E25" = "  +0.57%  "
    "D26" = "12.05"
    "E26" = "  -0.93%  "
    "D27" = "2.12"
    "E27" = "  -0.65%  "
    "E28" = "  -0.39%  "
    "E29" = "  +0.14%  "
    "D30" = "3.960.44"
    "E30" = "  -0.20%  "
    "D31" = "2.64"
    "E31" = "  -5.62%  "
    "B32" = "ImmutableX"
    "C32" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
    "D32" = "2.23"
    "E32" = "  -0.30%  "
    "B33" = "NEARProtocol"
    "C33" = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
    "D33" = "7.29"
    "E33" = "  -1.20%  "
    "D34" = "29.24"
    "E34" = "  -1.06%  "
    "D35" = "0.999"
    "E35" = "  -0.23%  "
    "D36" = "9.03"
    "E36" = "  -1.02%  "
    "D37" = "0.101"
    "E37" = "  +0.91%  "
    "D38" = "0.149"
    "E38" = "  +8.27%  "
    "D39" = "5.91"
    "E39" = "  +1.59%  "
    "D40" = "3.25"
    "E40" = "  -1.24%  "
    "D41" = "0.981"
    "E41" = "  -1.81%  "
    "D42" = "1.00"
    "E42" = "  +0.09%  "
    "E43" = "  +0.01%  "
    "B44" = "TheGraph"
    "C44" = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
    "D44" = "0.298"
    "E44" = "  -0.80%  "
    "D45" = "43.27"
    "E45" = "  -2.81%  "
    "D46" = "47.18"
    "E46" = "  -1.14%  "
    "B47" = "Monero"
    "C47" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    "D47" = "153.13"
    "E47" = "  +1.43%  "
    "B48" = "ONDO"
    "C48" = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
    "D48" = "1.40"
    "E48" = "  -1.85%  "
    "D49" = "8.40"
    "E49" = "  +0.41%  "
    "D50" = "1.87"
    "E50" = "  +0.66%  "
    "B51" = "Bittensor"
    "C51" = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
    "D51" = "384.05"
    "E51" = "  -1.42%  "
}

foreach ($ref in $updates.Keys) {
    $col = $ref -replace "[0-9]+$", ""
    $value = $updates[$ref]
    if ($col -eq "D") {
        # Column D holds numeric-looking strings (prices) that must remain text
        $ws.Range($ref).NumberFormat = "@"
        $ws.Range($ref).Value = $value
        $ws.Range($ref).Style = "Normal"
    } else {
        $ws.Range($ref).Value = $value
    }
}
